$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room: insert 22 new full rows at 19..40. This pushes the old
#    footer rows (23 "signature line", 24 "legal rep name / firm") down to
#    45/46 automatically (with their merged cells following).
# ---------------------------------------------------------------------------
$ws.Range("19:40").Insert()

# ---------------------------------------------------------------------------
# 2. Propagate formatting.
#    Row 18 still carries the old "last row" thick-bottom-border style;
#    copy it down to the new last row (40) before it gets overwritten.
#    Then stamp the plain interior-row style (from row 17) across rows
#    18..39 so every detail row matches the table's normal look.
# ---------------------------------------------------------------------------
$ws.Range("B18:J18").Copy($ws.Range("B40:J40"))
$ws.Range("B17:J17").Copy($ws.Range("B18:J39"))

# ---------------------------------------------------------------------------
# 3. Fill the worker / late-payment-period detail rows (16-40).
#    Columns: B=TipoDoc, C=NumDoc, D=Nombre, E=PeriodoMora, F=ValorMora,
#    G=SalarioBasico
# ---------------------------------------------------------------------------
$rows = @(
  @{r=16; b="CC"; c="1051448224"; d="CINDY ACEVEDO TILVEZ";               e="2506"; f=68000; g=1160000},
  @{r=17; b="CC"; c="1051448224"; d="CINDY ACEVEDO TILVEZ";               e="2505"; f=68000; g=1160000},
  @{r=18; b="CC"; c="1051448224"; d="CINDY ACEVEDO TILVEZ";               e="2504"; f=68000; g=1160000},
  @{r=19; b="CC"; c="1051448224"; d="CINDY ACEVEDO TILVEZ";               e="2503"; f=22667; g=1160000},
  @{r=20; b="CC"; c="1047372540"; d="JOSE MANUEL RODRIGUEZ SERRANO";      e="2507"; f=52000; g=877803},
  @{r=21; b="CC"; c="1047372540"; d="JOSE MANUEL RODRIGUEZ SERRANO";      e="2506"; f=52000; g=877803},
  @{r=22; b="CC"; c="1047372540"; d="JOSE MANUEL RODRIGUEZ SERRANO";      e="2505"; f=52000; g=877803},
  @{r=23; b="CC"; c="1047372540"; d="JOSE MANUEL RODRIGUEZ SERRANO";      e="2504"; f=52000; g=877803},
  @{r=24; b="CC"; c="1047372540"; d="JOSE MANUEL RODRIGUEZ SERRANO";      e="2503"; f=52000; g=877803},
  @{r=25; b="CC"; c="1047372540"; d="JOSE MANUEL RODRIGUEZ SERRANO";      e="2502"; f=52000; g=877803},
  @{r=26; b="CC"; c="1143328043"; d="CRISTIAN JOSE ROMERO ARIAS";         e="2507"; f=56940; g=1423500},
  @{r=27; b="CC"; c="1143328043"; d="CRISTIAN JOSE ROMERO ARIAS";         e="2506"; f=56940; g=1423500},
  @{r=28; b="CC"; c="1143328043"; d="CRISTIAN JOSE ROMERO ARIAS";         e="2505"; f=56940; g=1423500},
  @{r=29; b="CC"; c="1143328043"; d="CRISTIAN JOSE ROMERO ARIAS";         e="2504"; f=56940; g=1423500},
  @{r=30; b="CC"; c="1143328043"; d="CRISTIAN JOSE ROMERO ARIAS";         e="2503"; f=56940; g=1423500},
  @{r=31; b="CC"; c="1143328043"; d="CRISTIAN JOSE ROMERO ARIAS";         e="2502"; f=56940; g=1423500},
  @{r=32; b="CC"; c="1137195038"; d="JESUS DAVID CANTILLO GUERRERO";      e="2507"; f=80000; g=2000000},
  @{r=33; b="CC"; c="1137195038"; d="JESUS DAVID CANTILLO GUERRERO";      e="2506"; f=80000; g=2000000},
  @{r=34; b="CC"; c="1137195038"; d="JESUS DAVID CANTILLO GUERRERO";      e="2505"; f=80000; g=2000000},
  @{r=35; b="CC"; c="1137195038"; d="JESUS DAVID CANTILLO GUERRERO";      e="2504"; f=80000; g=2000000},
  @{r=36; b="CC"; c="1137195038"; d="JESUS DAVID CANTILLO GUERRERO";      e="2503"; f=80000; g=2000000},
  @{r=37; b="CC"; c="1137195038"; d="JESUS DAVID CANTILLO GUERRERO";      e="2502"; f=80000; g=2000000},
  @{r=38; b="CC"; c="1047447516"; d="SOL DALIS BLANCO CASTRO";            e="2507"; f=2400;  g=1700000},
  @{r=39; b="CC"; c="1051446826"; d="YURIS DEL CARMEN CASTRO PATERNINA";  e="2506"; f=13600; g=1700000},
  @{r=40; b="CC"; c="1002412621"; d="FERNI JOSE ACEVEDO GOMEZ";           e="2507"; f=10333; g=1550000}
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r, 2).Value = $row.b
  $ws.Cells.Item($r, 3).Value = $row.c
  $ws.Cells.Item($r, 4).Value = $row.d
  $ws.Cells.Item($r, 5).Value = $row.e
  $ws.Cells.Item($r, 6).Value = $row.f
  $ws.Cells.Item($r, 7).Value = $row.g
}

# ---------------------------------------------------------------------------
# 4. Update the summary header: total overdue value, worker count, period
#    count (labels themselves are unchanged).
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1386640
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 6
